# Generate Report for Handback
# - 7abf3104...md and c771cdfd...md are now out of sync with en-US (their
#   handback status flips from "in sync" to "not in sync").
# - c771cdfd...md just got handed back: stamp its "Correspond Handback
#   DateTime" for both the zh-cn and de-de locales.
# - Widen the status columns so the longer "not in sync" text still fits.

$wb = $excel.ActiveWorkbook

$inSync    = "Handed back: in sync with en-US"
$notInSync = "Handed back: not in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: columns E (zh-cn) / F (de-de) show the handback status
# for each source file. Rows: 2 = 7abf3104, 3 = c771cdfd, 4 = f22cedc0.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $notInSync
$wsOverview.Range("F2").Value = $notInSync
$wsOverview.Range("E3").Value = $notInSync
$wsOverview.Range("F3").Value = $notInSync

# ---------------------------------------------------------------------
# zh-cn sheet: column C = Status, column L = Correspond Handback DateTime
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $notInSync
$wsZhCn.Range("C3").Value = $notInSync
$wsZhCn.Range("L3").Value = "2017-01-03 04:23:04"

# ---------------------------------------------------------------------
# de-de sheet: column C = Status, column L = Correspond Handback DateTime
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $notInSync
$wsDeDe.Range("C3").Value = $notInSync
$wsDeDe.Range("L3").Value = "2017-01-03 04:23:15"

# ---------------------------------------------------------------------
# Widen the status columns to fit the longer "not in sync" text.
# ---------------------------------------------------------------------
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = 33.42
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 33.42
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 33.42
